$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 3977.875
$ws.Range("I34").Value = 689
$ws.Range("K34").Value = 689
$ws.Range("M34").Value = -486
$ws.Range("H36").Value = 3977.875
$ws.Range("I36").Value = 689
$ws.Range("K36").Value = 689
$ws.Range("M36").Value = 26
$ws.Range("H40").Value = 2665
$ws.Range("I40").Value = 5466.6665
$ws.Range("K40").Value = 5466.6665
$ws.Range("M40").Value = -5291.6665
$ws.Range("H53").Value = 165.93333
$ws.Range("I53").Value = 86.90000000000001
$ws.Range("J53").Value = 324
$ws.Range("K53").Value = 86.90000000000001
$ws.Range("L53").Value = 324
$ws.Range("M53").Value = 550.1
$ws.Range("N53").Value = -1598
$ws.Range("H62").Value = 5360.8237
$ws.Range("I62").Value = 3138.2144
$ws.Range("J62").Value = 15733
$ws.Range("K62").Value = 3138.2144
$ws.Range("L62").Value = 15733
$ws.Range("M62").Value = -2514.2144
$ws.Range("N62").Value = -16981
$ws.Range("H64").Value = 3859.3
$ws.Range("I64").Value = 3513.2856
$ws.Range("J64").Value = 4666.6665
$ws.Range("K64").Value = 3513.2856
$ws.Range("L64").Value = 4666.6665
$ws.Range("M64").Value = -3265.2856
$ws.Range("N64").Value = -5162.6665
$ws.Range("H65").Value = 5360.8237
$ws.Range("I65").Value = 3138.2144
$ws.Range("J65").Value = 15733
$ws.Range("K65").Value = 15691.072
$ws.Range("L65").Value = 78665
$ws.Range("M65").Value = -12571.072
$ws.Range("N65").Value = -84905
$ws.Range("H67").Value = 3859.3
$ws.Range("I67").Value = 3513.2856
$ws.Range("J67").Value = 4666.6665
$ws.Range("K67").Value = 3513.2856
$ws.Range("L67").Value = 4666.6665
$ws.Range("M67").Value = -2655.2856
$ws.Range("N67").Value = -6382.6665
$ws.Range("H70").Value = 2000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 2000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 6000
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -6540
$ws.Range("H73").Value = 2000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 2000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 6000
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -7872

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1292.4412
$ws.Range("I45").Value = 1313.7858
$ws.Range("K45").Value = 1313.7858
$ws.Range("M45").Value = -936.7858000000001
$ws.Range("H64").Value = 33318.25
$ws.Range("I64").Value = 27000
$ws.Range("J64").Value = 35424.332
$ws.Range("K64").Value = 27000
$ws.Range("L64").Value = 35424.332
$ws.Range("M64").Value = -26752
$ws.Range("N64").Value = -35920.332
$ws.Range("H67").Value = 33318.25
$ws.Range("I67").Value = 27000
$ws.Range("J67").Value = 35424.332
$ws.Range("K67").Value = 27000
$ws.Range("L67").Value = 35424.332
$ws.Range("M67").Value = -26142
$ws.Range("N67").Value = -37140.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 30590
$ws.Range("I62").Value = 22999
$ws.Range("J62").Value = 38181
$ws.Range("K62").Value = 22999
$ws.Range("L62").Value = 38181
$ws.Range("M62").Value = -22313
$ws.Range("N62").Value = -39553
$ws.Range("H65").Value = 30590
$ws.Range("I65").Value = 22999
$ws.Range("J65").Value = 38181
$ws.Range("K65").Value = 68997
$ws.Range("L65").Value = 114543
$ws.Range("M65").Value = -65565
$ws.Range("N65").Value = -121407

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 679.8
$ws.Range("I107").Value = 393.14285
$ws.Range("J107").Value = 834.1539
$ws.Range("K107").Value = 393.14285
$ws.Range("L107").Value = 834.1539
$ws.Range("M107").Value = 1526.85715
$ws.Range("N107").Value = -4674.1539

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9220.583000000001
$ws.Range("I3").Value = 3775
$ws.Range("J3").Value = 14666.167
$ws.Range("K3").Value = 11325
$ws.Range("L3").Value = 43998.501
$ws.Range("M3").Value = -11213
$ws.Range("N3").Value = -44222.501
$ws.Range("H114").Value = 452.36365
$ws.Range("I114").Value = 347.5
$ws.Range("J114").Value = 512.2857
$ws.Range("K114").Value = 1042.5
$ws.Range("L114").Value = 1536.8571
$ws.Range("M114").Value = 2211.5
$ws.Range("N114").Value = -8044.8571
$ws.Range("H131").Value = 852.13
$ws.Range("J131").Value = 896.2778
$ws.Range("L131").Value = 2688.8334
$ws.Range("N131").Value = -12768.8334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3551.7896
$ws.Range("I7").Value = 3725.0833
$ws.Range("J7").Value = 3254.7144
$ws.Range("K7").Value = 3725.0833
$ws.Range("L7").Value = 3254.7144
$ws.Range("M7").Value = -3613.0833
$ws.Range("N7").Value = -3478.7144
$ws.Range("H22").Value = 615.125
$ws.Range("I22").Value = 333.8
$ws.Range("K22").Value = 333.8
$ws.Range("M22").Value = -38.80000000000001
$ws.Range("H27").Value = 615.125
$ws.Range("I27").Value = 333.8
$ws.Range("K27").Value = 333.8
$ws.Range("M27").Value = -226.8
$ws.Range("H68").Value = 2475.5
$ws.Range("I68").Value = 2300.6667
$ws.Range("K68").Value = 2300.6667
$ws.Range("M68").Value = -1551.6667
$ws.Range("H71").Value = 2475.5
$ws.Range("I71").Value = 2300.6667
$ws.Range("K71").Value = 11503.3335
$ws.Range("M71").Value = -7759.333500000001
$ws.Range("H122").Value = 5802.769
$ws.Range("I122").Value = 5352.472
$ws.Range("K122").Value = 16057.416
$ws.Range("M122").Value = -13607.416
$ws.Range("H126").Value = 3551.7896
$ws.Range("I126").Value = 3725.0833
$ws.Range("J126").Value = 3254.7144
$ws.Range("K126").Value = 11175.2499
$ws.Range("L126").Value = 9764.143199999999
$ws.Range("M126").Value = -8705.249899999999
$ws.Range("N126").Value = -14704.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3751
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 3573.1428
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 3573.1428
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -4821.1428
$ws.Range("H65").Value = 3751
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 3573.1428
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 17865.714
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -24105.714

